$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.313.55'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '3.767.46'
$ws.Range("E3").Value = '  -1.59%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '693.11'
$ws.Range("D5").Style = 'Normal'
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '167.39'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -2.50%  '
$ws.Range("D7").Value = '3.766.49'
$ws.Range("E7").Value = '  -1.58%  '
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("E12").Value = '  +3.15%  '
$ws.Range("E13").Value = '  -3.49%  '
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").Value = '4.404.25'
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("D16").Value = '3.768.04'
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("D17").Value = '70.406.15'
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '17.31'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '512.06'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +3.48%  '
$ws.Range("E22").Value = '  -3.30%  '
$ws.Range("E23").Value = '  -3.93%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '83.17'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -2.52%  '
$ws.Range("E25").Value = '  -4.96%  '
$ws.Range("E26").Value = '  +3.10%  '
$ws.Range("D27").Value = '3.915.60'
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("E28").Value = '  -4.86%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -7.41%  '
$ws.Range("E31").Value = '  -6.25%  '
$ws.Range("E32").Value = '  -1.37%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '7.23'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '28.79'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -1.97%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("E36").Value = '  -5.28%  '
$ws.Range("E37").Value = '  -0.48%  '
$ws.Range("D38").Value = '3.732.87'
$ws.Range("E38").Value = '  -1.53%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '6.48'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +8.28%  '
$ws.Range("E40").Value = '  -3.48%  '
$ws.Range("E42").Value = '  -3.94%  '
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("E45").Value = '  -6.78%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '160.73'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -1.60%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '48.88'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").Value = '  -5.24%  '
$ws.Range("E49").Value = '  -2.81%  '
$ws.Range("E50").Value = '  -0.65%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '406.58'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -5.18%  '
